# Generate Report for Handoff
# Updates the handoff-priority column ("Priority") to "ht" and refreshes
# the related handoff/generate timestamps for the rows that were just
# (re-)handed off.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 13)

# "Overview" sheet: bump "Latest HO Xliff Generate Date" (column G)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-15 16:21:33"
}

# "zh-cn" sheet: mark Priority (column E) as "ht" and bump
# "Latest Handoff Datetime" (column H)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-15 16:21:28"
}

# "de-de" sheet: mark Priority (column E) as "ht" and bump
# "Latest Handoff Datetime" (column H)
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-15 16:21:33"
}
